$wb = $excel.ActiveWorkbook

# "final dataset" is the sheet that gets the new simulation/test columns.
$ws = $wb.Worksheets.Item("final dataset")

# Make it the active tab/sheet (matches activeTab="2" + tabSelected on this
# sheet in the saved workbook) BEFORE the structural edit below, so the
# selection set afterwards sticks to this sheet.
$ws.Activate() | Out-Null

# Remove the old "title of paper" column (column B) entirely - this shifts
# every later column one to the left (title of paper's shared string simply
# stops being referenced and drops out of the shared-string table).
$ws.Columns("B").Delete() | Out-Null

# The column that used to be "top_10_frequent_words_overall" (now column K
# after the shift) is renamed to reflect the new summarization width.
$ws.Range("K1").Value = "top_5_frequent_words_overall"

# Leave the selection where the author left it when they saved.
$ws.Range("K2").Select() | Out-Null
